$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Flip existing "Can exchange eligible" flags (H column) from 0 to 1 ---
$rowsToFlip = @(13,21,23,24,25,26,27,28,29)
foreach ($r in $rowsToFlip) {
    $ws.Cells.Item($r, 8).Value = 1
}

# --- 2) Add two new rows (30 and 31) for the new 2024 Spanish coins ---
# Copy formatting from row 29 down to rows 30 and 31 first.
$srcRow = $ws.Range("A29:I29")
$srcRow.Copy()
$ws.Range("A30:I30").PasteSpecial(-4122) | Out-Null
$srcRow.Copy()
$ws.Range("A31:I31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 30: 200th Anniversary - National Police Corps
$ws.Cells.Item(30, 1).Value = 2024
$ws.Cells.Item(30, 2).Value = "200th Anniversary - National Police Corps"
$ws.Cells.Item(30, 4).Value = "Obv: With mint logo - ""M"""
$ws.Cells.Item(30, 5).Value = "Rev: new map of Europe"
$ws.Cells.Item(30, 8).Value = 1
$ws.Range("I30").Formula = '=IF(OR(AND(H30>1,H30<>"-")),"Can exchange","")'

# Row 31: Sevilla (Spanish UNESCO series)
$ws.Cells.Item(31, 1).Value = 2024
$ws.Cells.Item(31, 2).Value = "Sevilla"
$ws.Cells.Item(31, 3).Value = "Spanish UNESCO"
$ws.Cells.Item(31, 4).Value = "Obv: With mint logo - ""M"""
$ws.Cells.Item(31, 5).Value = "Rev: new map of Europe"
$ws.Cells.Item(31, 8).Value = 1
$ws.Range("I31").Formula = '=IF(OR(AND(H31>1,H31<>"-")),"Can exchange","")'

# --- 3) Conditional formatting for the two new rows (mirrors the pattern used for each data row) ---
function Add-ExchangeConditionalFormatting($rangeAddress) {
    $rng = $ws.Range($rangeAddress)
    $fcText = $rng.FormatConditions.Add(9, 0, $null, $null, "*-")
    $fcText.Formula1 = "=NOT(ISERROR(SEARCH((" + [char]34 + "*-" + [char]34 + "),(" + $rangeAddress + "))))"

    $fcColor = $rng.FormatConditions.AddColorScale(3)
    $fcColor.ColorScaleCriteria.Item(1).Type = 0
    $fcColor.ColorScaleCriteria.Item(1).Value = 0
    $fcColor.ColorScaleCriteria.Item(1).FormatColor.Color = 10184191
    $fcColor.ColorScaleCriteria.Item(2).Type = 0
    $fcColor.ColorScaleCriteria.Item(2).Value = 1
    $fcColor.ColorScaleCriteria.Item(2).FormatColor.Color = 11719377
    $fcColor.ColorScaleCriteria.Item(3).Type = 0
    $fcColor.ColorScaleCriteria.Item(3).Value = 10
    $fcColor.ColorScaleCriteria.Item(3).FormatColor.Color = 5287680
}

Add-ExchangeConditionalFormatting("H30")
Add-ExchangeConditionalFormatting("H31")

# --- 4) Cosmetic: selection on sheet1 and sheet2 like the saved file ---
$ws.Range("J8").Select()
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B5").Select()
$ws.Select()
